$d = $word.ActiveDocument

# Locate the target empty paragraph (the one right after
# "Don't use delete topic, it will lead to lot of errors.") and the
# "Producing topic -" paragraph right after it.
$targetIndex = -1
$producingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $prevText = ""
    if ($i -gt 1) { $prevText = $d.Paragraphs($i - 1).Range.Text }
    if ($prevText -match "Don.t use delete topic" -and $p.Range.Text.Trim().Length -eq 0) {
        $targetIndex = $i
    }
    if ($targetIndex -ne -1 -and $i -eq ($targetIndex + 1) -and $p.Range.Text.Trim() -eq "Producing topic –") {
        $producingIndex = $i
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate target empty paragraph"
}

# Replace the "Producing topic -" paragraph FIRST (so that the earlier
# paragraph's index / range stays valid for the subsequent insert).
if ($producingIndex -ne -1) {
    $producingPara = $d.Paragraphs($producingIndex)
    $producingXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Producing topic –</w:t></w:r></w:p>'
    $producingPara.Range.InsertXML($producingXml)
}

$targetPara = $d.Paragraphs($targetIndex)
$targetRange = $targetPara.Range

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Note: use below commands in new version of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-US"/></w:rPr><w:t>kafka</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> –</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>create topic -</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">kafka-topics.bat --bootstrap-server localhost:9092 --topic </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>first_topic</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> --create --partitions 1 --replication-factor 1</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>list topic -</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">kafka-topics.bat --bootstrap-server localhost:9092 --list </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>Describe topic -</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">kafka-topics.bat --bootstrap-server localhost:9092 --topic </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>first_topic</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> –describe</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Delete topic -</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">kafka-topics.bat --bootstrap-server localhost:9092 --topic </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>first_topic</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> –delete</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetRange.InsertXML($xml)

Write-Host "Done"
